# Update "想去人数" (number of interested attendees) counts in column F
# across the workbook's four sheets, per the upstream data refresh
# ("Update gh-pages to output generated at 74db155").

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 253
$ws1.Range("F5").Value  = 1996
$ws1.Range("F6").Value  = 80
$ws1.Range("F7").Value  = 456
$ws1.Range("F8").Value  = 418
$ws1.Range("F10").Value = 7094
$ws1.Range("F12").Value = 548
$ws1.Range("F13").Value = 428
$ws1.Range("F14").Value = 70
$ws1.Range("F15").Value = 2415
$ws1.Range("F16").Value = 1775
$ws1.Range("F17").Value = 149
$ws1.Range("F18").Value = 49
$ws1.Range("F19").Value = 105
$ws1.Range("F21").Value = 118
$ws1.Range("F23").Value = 173
$ws1.Range("F24").Value = 82
$ws1.Range("F25").Value = 980
$ws1.Range("F26").Value = 177
$ws1.Range("F27").Value = 4118

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 13

# Sheet "本地生活" (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 126
$ws3.Range("F3").Value = 701

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 126
$ws4.Range("F4").Value  = 701
$ws4.Range("F7").Value  = 253
$ws4.Range("F8").Value  = 1996
$ws4.Range("F10").Value = 13
$ws4.Range("F11").Value = 80
$ws4.Range("F12").Value = 456
$ws4.Range("F13").Value = 418
$ws4.Range("F15").Value = 7094
$ws4.Range("F17").Value = 548
$ws4.Range("F18").Value = 428
$ws4.Range("F19").Value = 70
$ws4.Range("F20").Value = 2415
$ws4.Range("F21").Value = 1775
$ws4.Range("F22").Value = 149
$ws4.Range("F23").Value = 49
$ws4.Range("F24").Value = 105
$ws4.Range("F26").Value = 118
$ws4.Range("F28").Value = 173
$ws4.Range("F29").Value = 82
$ws4.Range("F30").Value = 980
$ws4.Range("F31").Value = 177
$ws4.Range("F32").Value = 4118
